$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Reorder the "completed_courses" values in row 3 (B3:I3)
$ws.Range("B3").Value = "CPSC 6179"
$ws.Range("C3").Value = "CYBR 6126"
$ws.Range("D3").Value = "CPSC 6185"
$ws.Range("E3").Value = "CPSC 6109"
$ws.Range("F3").Value = "CPSC 6177"
$ws.Range("G3").Value = "CPSC 6175"
$ws.Range("H3").Value = "CPSC 6127"
$ws.Range("I3").Value = "CPSC 6119"

# Update the error message in B7 to reflect the new ordering
$ws.Range("B7").Value = "Completed courses scheduled again: ['CPSC 6179', 'CYBR 6126', 'CPSC 6185', 'CPSC 6177', 'CPSC 6175', 'CPSC 6127', 'CPSC 6109', 'CPSC 6119']"
